# Adicionando as quantidades somadas dos meses de Janeiro, Fevereiro e Março,
# na planilha resumo
$wb = $excel.ActiveWorkbook

$meses = @("Janeiro", "Fevereiro", "Março")
$resumo = $wb.Worksheets.Item("Resumo")

# Para cada produto listado na coluna A da planilha Resumo (linhas 2 a 6),
# soma a Quantidade (coluna B) correspondente em cada planilha de mês.
for ($r = 2; $r -le 6; $r++) {
    $produto = $resumo.Cells.Item($r, 1).Value2
    $total = 0

    foreach ($mes in $meses) {
        $ws = $wb.Worksheets.Item($mes)
        for ($mr = 2; $mr -le 6; $mr++) {
            if ($ws.Cells.Item($mr, 1).Value2 -eq $produto) {
                $total = $total + $ws.Cells.Item($mr, 2).Value2
            }
        }
    }

    $resumo.Cells.Item($r, 2).Value = $total
}

$resumo.Select()
$resumo.Range("B2").Select()
